$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-04 18:38:22"

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
